# Update "Base ByteBank - Saldo do Cliente - Planilha Inicial.xlsx"
#
# The original sheet had a raw numeric ID (1001) in A2. The author replaced
# it with the text value "byte_1001" (introducing a new shared string) and
# left the selection sitting on A2 afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "byte_1001"

$ws.Range("A2").Select()
